# Auto-generated script applying cell value updates per the commit diff.
# Each worksheet is selected by name, then specific cells are updated via .Value,
# or cleared via .ClearContents() when the diff removes the cell entirely.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 2354.1667  # H40: 1571.16 -> 2354.1667
$ws.Cells.Item(40, 9).Value = 2164.2856  # I40: 1562.7858 -> 2164.2856
$ws.Cells.Item(40, 10).Value = 2620  # J40: 1581.8182 -> 2620
$ws.Cells.Item(40, 11).Value = 2164.2856  # K40: 1562.7858 -> 2164.2856
$ws.Cells.Item(40, 12).Value = 2620  # L40: 1581.8182 -> 2620
$ws.Cells.Item(40, 13).Value = -1989.2856  # M40: -1387.7858 -> -1989.2856
$ws.Cells.Item(40, 14).Value = -2970  # N40: -1931.8182 -> -2970
$ws.Cells.Item(58, 8).Value = 910080.9399999999  # H58: 1527.5714 -> 910080.9399999999
$ws.Cells.Item(58, 9).Value = 312.16666  # I58: 448.25 -> 312.16666
$ws.Cells.Item(58, 10).Value = 2001803.4  # J58: 2966.6667 -> 2001803.4
$ws.Cells.Item(58, 11).Value = 936.4999799999999  # K58: 1344.75 -> 936.4999799999999
$ws.Cells.Item(58, 12).Value = 6005410.199999999  # L58: 8900.000100000001 -> 6005410.199999999
$ws.Cells.Item(58, 13).Value = -786.4999799999999  # M58: -1194.75 -> -786.4999799999999
$ws.Cells.Item(58, 14).Value = -6005710.199999999  # N58: -9200.000100000001 -> -6005710.199999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 100001  # H4: 33350.332 -> 100001
$ws.Cells.Item(4, 9).Value = 100001  # I4: 50005.5 -> 100001
$ws.Cells.Item(4, 10).Value = 0  # J4: 40 -> 0
$ws.Cells.Item(4, 11).Value = 100001  # K4: 50005.5 -> 100001
$ws.Cells.Item(4, 12).Value = 0  # L4: 40 -> 0
$ws.Cells.Item(4, 13).Value = -99885  # M4: -49889.5 -> -99885
$ws.Cells.Item(4, 14).ClearContents()  # N4: -272 -> (removed)
$ws.Cells.Item(5, 8).Value = 75051  # H5: 42906.285 -> 75051
$ws.Cells.Item(5, 9).Value = 66734  # I5: 33390.332 -> 66734
$ws.Cells.Item(5, 11).Value = 66734  # K5: 33390.332 -> 66734
$ws.Cells.Item(5, 13).Value = -66622  # M5: -33278.332 -> -66622
$ws.Cells.Item(6, 8).Value = 35167.332  # H6: 50002 -> 35167.332
$ws.Cells.Item(6, 9).Value = 33668  # I6: 50002 -> 33668
$ws.Cells.Item(6, 10).Value = 36666.668  # J6: 0 -> 36666.668
$ws.Cells.Item(6, 11).Value = 33668  # K6: 50002 -> 33668
$ws.Cells.Item(6, 12).Value = 36666.668  # L6: 0 -> 36666.668
$ws.Cells.Item(6, 13).Value = -33495  # M6: -49829 -> -33495
$ws.Cells.Item(6, 14).Value = -37012.668  # N6: None -> -37012.668
$ws.Cells.Item(63, 8).Value = 9974.154  # H63: 10392.083 -> 9974.154
$ws.Cells.Item(63, 9).Value = 10597  # I63: 12040.5 -> 10597
$ws.Cells.Item(63, 10).Value = 2500  # J63: 2150 -> 2500
$ws.Cells.Item(63, 11).Value = 10597  # K63: 12040.5 -> 10597
$ws.Cells.Item(63, 12).Value = 2500  # L63: 2150 -> 2500
$ws.Cells.Item(63, 13).Value = -9911  # M63: -11354.5 -> -9911
$ws.Cells.Item(63, 14).Value = -3872  # N63: -3522 -> -3872
$ws.Cells.Item(66, 8).Value = 9974.154  # H66: 10392.083 -> 9974.154
$ws.Cells.Item(66, 9).Value = 10597  # I66: 12040.5 -> 10597
$ws.Cells.Item(66, 10).Value = 2500  # J66: 2150 -> 2500
$ws.Cells.Item(66, 11).Value = 52985  # K66: 60202.5 -> 52985
$ws.Cells.Item(66, 12).Value = 12500  # L66: 10750 -> 12500
$ws.Cells.Item(66, 13).Value = -49553  # M66: -56770.5 -> -49553
$ws.Cells.Item(66, 14).Value = -19364  # N66: -17614 -> -19364

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 75051  # H4: 42906.285 -> 75051
$ws.Cells.Item(4, 9).Value = 66734  # I4: 33390.332 -> 66734
$ws.Cells.Item(4, 11).Value = 66734  # K4: 33390.332 -> 66734
$ws.Cells.Item(4, 13).Value = -66619  # M4: -33275.332 -> -66619
$ws.Cells.Item(15, 8).Value = 19949.5  # H15: 33222 -> 19949.5
$ws.Cells.Item(15, 10).Value = 19949.5  # J15: 33222 -> 19949.5
$ws.Cells.Item(15, 12).Value = 19949.5  # L15: 33222 -> 19949.5
$ws.Cells.Item(15, 14).Value = -20403.5  # N15: -33676 -> -20403.5
$ws.Cells.Item(19, 8).Value = 40001500  # H19: 80000000 -> 40001500
$ws.Cells.Item(19, 10).Value = 40001500  # J19: 80000000 -> 40001500
$ws.Cells.Item(19, 12).Value = 40001500  # L19: 80000000 -> 40001500
$ws.Cells.Item(19, 14).Value = -40001846  # N19: -80000346 -> -40001846
$ws.Cells.Item(35, 8).Value = 43060  # H35: 42242.855 -> 43060
$ws.Cells.Item(35, 9).Value = 43500  # I35: 0 -> 43500
$ws.Cells.Item(35, 10).Value = 42950  # J35: 42242.855 -> 42950
$ws.Cells.Item(35, 11).Value = 43500  # K35: 0 -> 43500
$ws.Cells.Item(35, 12).Value = 42950  # L35: 42242.855 -> 42950
$ws.Cells.Item(35, 13).Value = -43190  # M35: None -> -43190
$ws.Cells.Item(35, 14).Value = -43570  # N35: -42862.855 -> -43570
$ws.Cells.Item(82, 8).Value = 0  # H82: 29250 -> 0
$ws.Cells.Item(82, 9).Value = 0  # I82: 3000 -> 0
$ws.Cells.Item(82, 10).Value = 0  # J82: 38000 -> 0
$ws.Cells.Item(82, 11).Value = 0  # K82: 3000 -> 0
$ws.Cells.Item(82, 12).Value = 0  # L82: 38000 -> 0
$ws.Cells.Item(82, 13).ClearContents()  # M82: -2617 -> (removed)
$ws.Cells.Item(82, 14).ClearContents()  # N82: -38766 -> (removed)
$ws.Cells.Item(85, 8).Value = 0  # H85: 29250 -> 0
$ws.Cells.Item(85, 9).Value = 0  # I85: 3000 -> 0
$ws.Cells.Item(85, 10).Value = 0  # J85: 38000 -> 0
$ws.Cells.Item(85, 11).Value = 0  # K85: 3000 -> 0
$ws.Cells.Item(85, 12).Value = 0  # L85: 38000 -> 0
$ws.Cells.Item(85, 13).ClearContents()  # M85: -1674 -> (removed)
$ws.Cells.Item(85, 14).ClearContents()  # N85: -40652 -> (removed)

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(17, 8).Value = 0  # H17: 2000 -> 0
$ws.Cells.Item(17, 9).Value = 0  # I17: 2000 -> 0
$ws.Cells.Item(17, 11).Value = 0  # K17: 2000 -> 0
$ws.Cells.Item(17, 13).ClearContents()  # M17: -1826 -> (removed)
$ws.Cells.Item(25, 8).Value = 2224  # H25: 1806.8 -> 2224
$ws.Cells.Item(25, 9).Value = 555  # I25: 527.5 -> 555
$ws.Cells.Item(25, 10).Value = 2502.1667  # J25: 2126.625 -> 2502.1667
$ws.Cells.Item(25, 11).Value = 555  # K25: 527.5 -> 555
$ws.Cells.Item(25, 12).Value = 2502.1667  # L25: 2126.625 -> 2502.1667
$ws.Cells.Item(25, 13).Value = -381  # M25: -353.5 -> -381
$ws.Cells.Item(25, 14).Value = -2850.1667  # N25: -2474.625 -> -2850.1667
$ws.Cells.Item(41, 8).Value = 7500  # H41: 0 -> 7500
$ws.Cells.Item(41, 10).Value = 7500  # J41: 0 -> 7500
$ws.Cells.Item(41, 12).Value = 7500  # L41: 0 -> 7500
$ws.Cells.Item(41, 14).Value = -8356  # N41: None -> -8356
$ws.Cells.Item(50, 8).Value = 41666.668  # H50: 40675 -> 41666.668
$ws.Cells.Item(50, 10).Value = 41666.668  # J50: 40675 -> 41666.668
$ws.Cells.Item(50, 12).Value = 41666.668  # L50: 40675 -> 41666.668
$ws.Cells.Item(50, 14).Value = -42916.668  # N50: -41925 -> -42916.668
$ws.Cells.Item(51, 8).Value = 38311.11  # H51: 37072.727 -> 38311.11
$ws.Cells.Item(51, 10).Value = 38311.11  # J51: 37072.727 -> 38311.11
$ws.Cells.Item(51, 12).Value = 38311.11  # L51: 37072.727 -> 38311.11
$ws.Cells.Item(51, 14).Value = -39783.11  # N51: -38544.727 -> -39783.11
$ws.Cells.Item(59, 8).Value = 42250  # H59: 41777.777 -> 42250
$ws.Cells.Item(59, 10).Value = 42250  # J59: 41777.777 -> 42250
$ws.Cells.Item(59, 12).Value = 42250  # L59: 41777.777 -> 42250
$ws.Cells.Item(59, 14).Value = -44540  # N59: -44067.777 -> -44540
$ws.Cells.Item(60, 8).Value = 14778.28  # H60: 15693.228 -> 14778.28
$ws.Cells.Item(60, 9).Value = 11000  # I60: 20000 -> 11000
$ws.Cells.Item(60, 10).Value = 15106.826  # J60: 15488.143 -> 15106.826
$ws.Cells.Item(60, 11).Value = 11000  # K60: 20000 -> 11000
$ws.Cells.Item(60, 12).Value = 15106.826  # L60: 15488.143 -> 15106.826
$ws.Cells.Item(60, 13).Value = -10489  # M60: -19489 -> -10489
$ws.Cells.Item(60, 14).Value = -16128.826  # N60: -16510.143 -> -16128.826
$ws.Cells.Item(61, 8).Value = 38311.11  # H61: 37072.727 -> 38311.11
$ws.Cells.Item(61, 10).Value = 38311.11  # J61: 37072.727 -> 38311.11
$ws.Cells.Item(61, 12).Value = 38311.11  # L61: 37072.727 -> 38311.11
$ws.Cells.Item(61, 14).Value = -39007.11  # N61: -37768.727 -> -39007.11
$ws.Cells.Item(68, 8).Value = 10000  # H68: 30295 -> 10000
$ws.Cells.Item(68, 9).Value = 10000  # I68: 0 -> 10000
$ws.Cells.Item(68, 10).Value = 0  # J68: 30295 -> 0
$ws.Cells.Item(68, 11).Value = 10000  # K68: 0 -> 10000
$ws.Cells.Item(68, 12).Value = 0  # L68: 30295 -> 0
$ws.Cells.Item(68, 13).Value = -9251  # M68: None -> -9251
$ws.Cells.Item(68, 14).ClearContents()  # N68: -31793 -> (removed)
$ws.Cells.Item(71, 8).Value = 10000  # H71: 30295 -> 10000
$ws.Cells.Item(71, 9).Value = 10000  # I71: 0 -> 10000
$ws.Cells.Item(71, 10).Value = 0  # J71: 30295 -> 0
$ws.Cells.Item(71, 11).Value = 30000  # K71: 0 -> 30000
$ws.Cells.Item(71, 12).Value = 0  # L71: 90885 -> 0
$ws.Cells.Item(71, 13).Value = -26256  # M71: None -> -26256
$ws.Cells.Item(71, 14).ClearContents()  # N71: -98373 -> (removed)
$ws.Cells.Item(74, 8).Value = 19814  # H74: 0 -> 19814
$ws.Cells.Item(74, 10).Value = 19814  # J74: 0 -> 19814
$ws.Cells.Item(74, 12).Value = 19814  # L74: 0 -> 19814
$ws.Cells.Item(74, 14).Value = -21562  # N74: None -> -21562
$ws.Cells.Item(77, 8).Value = 19814  # H77: 0 -> 19814
$ws.Cells.Item(77, 10).Value = 19814  # J77: 0 -> 19814
$ws.Cells.Item(77, 12).Value = 59442  # L77: 0 -> 59442
$ws.Cells.Item(77, 14).Value = -68178  # N77: None -> -68178
$ws.Cells.Item(94, 8).Value = 1529.0952  # H94: 1471.1818 -> 1529.0952
$ws.Cells.Item(94, 10).Value = 1260.2142  # J94: 1193.2 -> 1260.2142
$ws.Cells.Item(94, 12).Value = 1260.2142  # L94: 1193.2 -> 1260.2142
$ws.Cells.Item(94, 14).Value = -2162.2142  # N94: -2095.2 -> -2162.2142

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 1281.0857  # H113: 1190.1794 -> 1281.0857
$ws.Cells.Item(113, 9).Value = 1729.3334  # I113: 1506.6111 -> 1729.3334
$ws.Cells.Item(113, 10).Value = 944.9  # J113: 918.9524 -> 944.9
$ws.Cells.Item(113, 11).Value = 5188.0002  # K113: 4519.8333 -> 5188.0002
$ws.Cells.Item(113, 12).Value = 2834.7  # L113: 2756.8572 -> 2834.7
$ws.Cells.Item(113, 13).Value = -3018.0002  # M113: -2349.8333 -> -3018.0002
$ws.Cells.Item(113, 14).Value = -7174.7  # N113: -7096.8572 -> -7174.7
$ws.Cells.Item(131, 8).Value = 2395.3066  # H131: 2375.6052 -> 2395.3066
$ws.Cells.Item(131, 10).Value = 2700.7385  # J131: 2673.4243 -> 2700.7385
$ws.Cells.Item(131, 12).Value = 8102.2155  # L131: 8020.2729 -> 8102.2155
$ws.Cells.Item(131, 14).Value = -18182.2155  # N131: -18100.2729 -> -18182.2155

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 4509.9  # H80: 3999.5 -> 4509.9
$ws.Cells.Item(80, 9).Value = 7125  # I80: 0 -> 7125
$ws.Cells.Item(80, 10).Value = 2766.5  # J80: 3999.5 -> 2766.5
$ws.Cells.Item(80, 11).Value = 7125  # K80: 0 -> 7125
$ws.Cells.Item(80, 12).Value = 2766.5  # L80: 3999.5 -> 2766.5
$ws.Cells.Item(80, 13).Value = -6127  # M80: None -> -6127
$ws.Cells.Item(80, 14).Value = -4762.5  # N80: -5995.5 -> -4762.5
$ws.Cells.Item(83, 8).Value = 4509.9  # H83: 3999.5 -> 4509.9
$ws.Cells.Item(83, 9).Value = 7125  # I83: 0 -> 7125
$ws.Cells.Item(83, 10).Value = 2766.5  # J83: 3999.5 -> 2766.5
$ws.Cells.Item(83, 11).Value = 35625  # K83: 0 -> 35625
$ws.Cells.Item(83, 12).Value = 13832.5  # L83: 19997.5 -> 13832.5
$ws.Cells.Item(83, 13).Value = -30633  # M83: None -> -30633
$ws.Cells.Item(83, 14).Value = -23816.5  # N83: -29981.5 -> -23816.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(18, 8).Value = 50000  # H18: 0 -> 50000
$ws.Cells.Item(18, 10).Value = 50000  # J18: 0 -> 50000
$ws.Cells.Item(18, 12).Value = 50000  # L18: 0 -> 50000
$ws.Cells.Item(18, 14).Value = -50344  # N18: None -> -50344
$ws.Cells.Item(20, 8).Value = 20000  # H20: 19833.334 -> 20000
$ws.Cells.Item(20, 9).Value = 0  # I20: 20000 -> 0
$ws.Cells.Item(20, 10).Value = 20000  # J20: 19800 -> 20000
$ws.Cells.Item(20, 11).Value = 0  # K20: 20000 -> 0
$ws.Cells.Item(20, 12).Value = 20000  # L20: 19800 -> 20000
$ws.Cells.Item(20, 13).ClearContents()  # M20: -19774 -> (removed)
$ws.Cells.Item(20, 14).Value = -20452  # N20: -20252 -> -20452
$ws.Cells.Item(61, 8).Value = 1471.4762  # H61: 2222.2222 -> 1471.4762
$ws.Cells.Item(61, 9).Value = 1327.8334  # I61: 2028.5714 -> 1327.8334
$ws.Cells.Item(61, 10).Value = 2333.3333  # J61: 2900 -> 2333.3333
$ws.Cells.Item(61, 11).Value = 1327.8334  # K61: 2028.5714 -> 1327.8334
$ws.Cells.Item(61, 12).Value = 2333.3333  # L61: 2900 -> 2333.3333
$ws.Cells.Item(61, 13).Value = -1125.8334  # M61: -1826.5714 -> -1125.8334
$ws.Cells.Item(61, 14).Value = -2737.3333  # N61: -3304 -> -2737.3333
$ws.Cells.Item(113, 8).Value = 1471.4762  # H113: 2222.2222 -> 1471.4762
$ws.Cells.Item(113, 9).Value = 1327.8334  # I113: 2028.5714 -> 1327.8334
$ws.Cells.Item(113, 10).Value = 2333.3333  # J113: 2900 -> 2333.3333
$ws.Cells.Item(113, 11).Value = 1327.8334  # K113: 2028.5714 -> 1327.8334
$ws.Cells.Item(113, 12).Value = 2333.3333  # L113: 2900 -> 2333.3333
$ws.Cells.Item(113, 13).Value = 842.1666  # M113: 141.4286 -> 842.1666
$ws.Cells.Item(113, 14).Value = -6673.3333  # N113: -7240 -> -6673.3333

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 33812.645  # H122: 36079.207 -> 33812.645
$ws.Cells.Item(122, 9).Value = 40103.348  # I122: 41671.68 -> 40103.348
$ws.Cells.Item(122, 10).Value = 1101  # J122: 1126.25 -> 1101
$ws.Cells.Item(122, 11).Value = 120310.044  # K122: 125015.04 -> 120310.044
$ws.Cells.Item(122, 12).Value = 3303  # L122: 3378.75 -> 3303
$ws.Cells.Item(122, 13).Value = -117860.044  # M122: -122565.04 -> -117860.044
$ws.Cells.Item(122, 14).Value = -8203  # N122: -8278.75 -> -8203
